# Insert a new weekly price record as row 47 on the active sheet.
# All existing rows from 47 downward shift down by one (old row 47 -> 48,
# ..., old row 79 -> 80), matching dimension A1:T79 -> A1:T80.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(47).EntireRow.Insert()

$ws.Range("A47").Value = 9
$ws.Range("B47").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C47").Value = "Metropolitana"
$ws.Range("D47").Value = 45090
$ws.Range("E47").Value = 13
$ws.Range("F47").Value = "Fruta"
$ws.Range("G47").Value = 100107
$ws.Range("H47").Value = "Otros"
$ws.Range("I47").Value = 100107001
$ws.Range("J47").Value = "Caqui"
$ws.Range("K47").Value = "Fuyu"
$ws.Range("L47").Value = "Primera"
$ws.Range("M47").Value = 470
$ws.Range("N47").Value = 10500
$ws.Range("O47").Value = 11000
$ws.Range("P47").Value = 10734
$ws.Range("Q47").Value = "`$/caja 15 kilos granel"
$ws.Range("R47").Value = "Provincia de Curicó"
$ws.Range("S47").Value = 716
$ws.Range("T47").Value = 15
